$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("g10.2")

$data = @{
    2  = @(-3.377448218355161, 2.15496247269007, 7.434495881238234)
    3  = @(0.2448820992536849, -0.1721910504140167, -1.008236814482644)
    4  = @(5.624925926329238, 0.6459332894185987, 6.267294698459325)
    5  = @(6.16826137189288, -6.935153294182605, 9.755477202390939)
    6  = @(-1.494625744833378, -6.157851153201799, 8.009006717074318)
    7  = @(-0.2554344193826941, -4.932397146504464, 2.809173804671983)
    8  = @(-0.8961894651313584, -4.264270422584227, 0.1319837197746532)
    9  = @(4.683070112298715, -1.149109367878443, 11.26013612946959)
    10 = @(-10.69158489251918, -5.600889391489416, -5.900781225340602)
    11 = @(-6.373529693431146, 9.024371965595002, -7.350727472305019)
    12 = @(0.04395139504043133, 7.481601158193651, -10.82795411725803)
    13 = @(-2.445450002465022, 2.990192558263849, -3.269703462328233)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 2).Value = $values[0]
    $ws.Cells.Item($row, 3).Value = $values[1]
    $ws.Cells.Item($row, 4).Value = $values[2]
}
